$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.1559305616838552
$ws.Range("D2").Value = 0.1559304748967024

$ws.Range("C3").Value = 0.1637229051272682
$ws.Range("D3").Value = 0.1637228384233602

$ws.Range("C4").Value = 0.08546269601881656
$ws.Range("D4").Value = 0.0854627095093803

$ws.Range("C5").Value = 0.07238432245797
$ws.Range("D5").Value = 0.07238423279098104

$ws.Range("C6").Value = 0.130098276253158
$ws.Range("D6").Value = 0.130098355694125

$ws.Range("C7").Value = 0.1839516980365401
$ws.Range("D7").Value = 0.1839517227747055

$ws.Range("C8").Value = 0.208449540422392
$ws.Range("D8").Value = 0.2084496659107455
